$d = $word.ActiveDocument

# Locate the "KEY ACHIEVEMENTS AND IMPACT" section so the edits below only
# touch bullets inside it (some of the old bullet text is duplicated
# verbatim elsewhere in the document, e.g. under "Partner - Siege
# Analytics", and must be left alone).
$startPara = $null
$endPara = $null
foreach ($p in @($d.Paragraphs)) {
    if ($startPara -eq $null -and $p.Range.Text -like "*KEY ACHIEVEMENTS AND IMPACT*") {
        $startPara = $p
    } elseif ($startPara -ne $null -and $p.Range.Text -like "*TECHNICAL SKILLS*") {
        $endPara = $p
        break
    }
}
$sectionStart = $startPara.Range.Start
$sectionEnd = $endPara.Range.Start

# 1) Replace the first four "Key Achievements" bullet paragraphs in place
#    with new impact-focused accomplishment statements.
$r = $d.Range($sectionStart, $sectionEnd)
$r.Find.Execute(
    "• Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%",
    $true, $false, $false, $false, $false, $true, 0, $false,
    "• Predictive excellence: Achieved 87% voter turnout accuracy vs. 71% industry standard",
    2) | Out-Null

$r = $d.Range($sectionStart, $sectionEnd)
$r.Find.Execute(
    "• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from ±4.2% to ±2.1%",
    $true, $false, $false, $false, $false, $true, 0, $false,
    "• Reduced polling margins from ±4.2% to ±2.1%",
    2) | Out-Null

$r = $d.Range($sectionStart, $sectionEnd)
$r.Find.Execute(
    "• Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations",
    $true, $false, $false, $false, $false, $true, 0, $false,
    "• Executive authority: Briefed Presidents, Congressmen, Senators, Governors on election integrity, voter sentiment and postmortem analysis",
    2) | Out-Null

$r = $d.Range($sectionStart, $sectionEnd)
$r.Find.Execute(
    "• Developed longitudinal data analysis methods using geospatial techniques that improved segmentation accuracy by 34% and survey incidence rates by 28%, reducing polling costs while increasing response quality",
    $true, $false, $false, $false, $false, $true, 0, $false,
    "• Methodological advancement: Improved segmentation accuracy 34% and survey incidence 28%",
    2) | Out-Null

# 2) Remove the two trailing "Key Achievements" bullet paragraphs that are
#    no longer wanted. Collect the indices first, then delete from the
#    highest index down so earlier indices stay valid as the document
#    shrinks.
$indices = @()
$i = 1
foreach ($p in @($d.Paragraphs)) {
    $t = $p.Range.Text
    if ($t -like "*Provided expert testimony and press briefings on electoral data integrity and demographic modeling accuracy*" -or
        $t -like "*Demystified FEC compliance through real-time processing systems enabling transparent campaign finance monitoring*") {
        $indices += $i
    }
    $i = $i + 1
}
$indices = $indices[($indices.Count - 1)..0]
foreach ($idx in $indices) {
    $d.Paragraphs.Item($idx).Range.Delete()
}
